$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a blank row below the existing row 5 (between two plain trade
# rows) so the copied formatting is the plain data-row style rather than
# the bold header row above row 5. Existing rows 6-18 shift down to 7-19;
# row 5 (the old newest trade) stays put for now.
$ws.Rows.Item(6).Insert()

# Move the old row-5 trade (date 46059, ...) down into the freshly
# inserted row 6.
$ws.Cells.Item(6, 1).Value = $ws.Cells.Item(5, 1).Value2
$ws.Cells.Item(6, 2).Value = $ws.Cells.Item(5, 2).Value2
$ws.Cells.Item(6, 3).Value = $ws.Cells.Item(5, 3).Value2
$ws.Cells.Item(6, 4).Value = $ws.Cells.Item(5, 4).Value2
$ws.Cells.Item(6, 5).Value = $ws.Cells.Item(5, 5).Value2
$ws.Cells.Item(6, 6).Value = $ws.Cells.Item(5, 6).Value2
$ws.Cells.Item(6, 7).Value = $ws.Cells.Item(5, 7).Value2
$ws.Cells.Item(6, 10).Formula = "=Index!`$C`$2"

# Overwrite row 5 with the newest trade entry.
$ws.Cells.Item(5, 1).Value = 46062
$ws.Cells.Item(5, 2).Value = "NSE"
$ws.Cells.Item(5, 3).Value = "Buy"
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 2053
$ws.Cells.Item(5, 6).Value = 4135.11
$ws.Cells.Item(5, 7).Value = "CN#252611665409"
$ws.Cells.Item(5, 8).Value = 4.11
$ws.Cells.Item(5, 9).Value = 25
$ws.Cells.Item(5, 10).Formula = "=Index!`$C`$2"
